# Updated symbol list on Sat Dec 17 03:30:35 UTC 2022 with GitHub Actions
# Re-applies refreshed Price (column D) and Volume(1h) (column E) values
# for the crypto tracker sheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay TEXT (e.g. trailing
# zeros like "3.390" or "0.005000" would be lost if stored as real numbers),
# exactly like the inline-string cells already in the workbook. A leading
# apostrophe forces Excel to keep the entry as text instead of coercing it
# to a number.

$ws.Range("D2").Value = "'228.58"
$ws.Range("D3").Value = "'22.45"
$ws.Range("D4").Value = "'5.323"
$ws.Range("D5").Value = "'0.05535"
$ws.Range("D6").Value = "'3.390"
$ws.Range("D7").Value = "'6.469"
$ws.Range("D8").Value = "'1.070"
$ws.Range("D9").Value = "'0.7741"
$ws.Range("D10").Value = "'0.1384"
$ws.Range("D11").Value = "'0.07423"
$ws.Range("D12").Value = "'0.03138"
$ws.Range("D13").Value = "'0.02940"
$ws.Range("D14").Value = "'0.09256"
$ws.Range("D15").Value = "'0.001659"
$ws.Range("D16").Value = "'3.256"
$ws.Range("D17").Value = "'0.04768"
$ws.Range("D18").Value = "'0.0005890"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006218"
$ws.Range("D20").Value = "'0.005231"
$ws.Range("D21").Value = "'0.001063"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.868"
$ws.Range("D26").Value = "'0.1282"
$ws.Range("D27").Value = "'0.0005000"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Value = "'0.03948"
$ws.Range("D41").Value = "'0.007124"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("D43").Value = "'0.1033"
$ws.Range("D44").Value = "'0.008556"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D47").Value = "'0.7852"
$ws.Range("D48").Value = "'0.04057"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.01010"
